# Weekly fruit/vegetable price update: a new daily record is inserted at
# row 334 (pushing the existing rows 334-360 down to 335-361), and the new
# row 334 is populated with the latest observation for Acelga at the Vega
# Modelo de Temuco market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 334, shifting rows 334:360 down to 335:361.
$ws.Range("A334").EntireRow.Insert()

# Populate the new row 334 with the new observation.
$ws.Range("A334").Value = 10
$ws.Range("B334").Value = "Vega Modelo de Temuco"
$ws.Range("C334").Value = "La Araucanía"
$ws.Range("D334").Value = 44783
$ws.Range("E334").Value = 9
$ws.Range("F334").Value = 100112009
$ws.Range("G334").Value = "Acelga"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 40
$ws.Range("K334").Value = 9000
$ws.Range("L334").Value = 9000
$ws.Range("M334").Value = 9000
$ws.Range("N334").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O334").Value = "Región Metropolitana"
$ws.Range("P334").Value = 750
$ws.Range("Q334").Value = 12
$ws.Range("R334").Value = "Hortaliza"
